# Mark item "c." ("En clase pedido no se usaron atributos de asociación
# ERROR, se definieron dos constructores ERROR (solo constructor
# completo).") as resolved by applying strikethrough formatting to the
# whole paragraph (label, spacer run and text run), matching the
# strikethrough already used on other corrected items in this document.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "En clase pedido no se usaron atributos de asociación ERROR",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -and $rng.Find.Found) {
    $para = $rng.Paragraphs(1)
    $para.Range.Font.StrikeThrough = $true
} else {
    throw "Target paragraph not found"
}
